# Applies the "Updated cryptos list" data refresh to Sheet1.
# For numeric-looking text in column D, force text formatting first so
# Excel doesn't silently coerce values like '6.099' into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.059.11'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '1.907.11'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7486'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.62'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3091'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.48'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -4.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06962'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08092'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7680'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('D13').Value = '1.902.32'
$ws.Range('E13').Value = '  -1.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.279'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.84'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.40%  '
$ws.Range('D16').Value = '30.073.46'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.099'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.14'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '239.84'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.88%  '
$ws.Range('D21').Value = '2.171.71'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.107'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.46%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.345'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.96'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.96'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1276'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.054'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -6.44%  '
$ws.Range('E30').Value = '  -1.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.533'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.320'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.074'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05380'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.303'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7441'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.719'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01968'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.304'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4479'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.14'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.972'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8354'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.688'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.91'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.871'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.21%  '
$ws.Range('D49').Value = '2.065.29'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.52'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.90%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4139'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.84%  '

Write-Host "Applied 97 cell updates"
